$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("db")

# Row 4: Test
$ws.Range("A4").Value = "Test"
$ws.Range("B4").Value = 10
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 68
$ws.Range("E4").Value = 25

# Row 5: yeet
$ws.Range("A5").Value = "yeet"
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 0

# Row 6: wwwww
$ws.Range("A6").Value = "wwwww"
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0
